$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "28.877.91"
$ws.Cells.Item(2,5).Value = "  +8.00%  "
$ws.Cells.Item(3,4).Value = "1.811.10"
$ws.Cells.Item(3,5).Value = "  +4.97%  "
$ws.Cells.Item(4,4).Value = "'0.9993"
$ws.Cells.Item(4,5).Value = "  +0.22%  "
$ws.Cells.Item(5,4).Value = "'248.02"
$ws.Cells.Item(5,5).Value = "  +3.16%  "
$ws.Cells.Item(6,4).Value = "'0.9996"
$ws.Cells.Item(6,5).Value = "  +0.14%  "
$ws.Cells.Item(7,4).Value = "'0.4950"
$ws.Cells.Item(7,5).Value = "  +2.57%  "
$ws.Cells.Item(8,4).Value = "'43.25"
$ws.Cells.Item(8,5).Value = "  +5.03%  "
$ws.Cells.Item(9,4).Value = "'0.2782"
$ws.Cells.Item(9,5).Value = "  +7.81%  "
$ws.Cells.Item(10,4).Value = "'0.06421"
$ws.Cells.Item(11,4).Value = "1.808.27"
$ws.Cells.Item(11,5).Value = "  +4.93%  "
$ws.Cells.Item(12,4).Value = "'16.81"
$ws.Cells.Item(12,5).Value = "  +5.92%  "
$ws.Cells.Item(13,5).Value = "  +3.71%  "
$ws.Cells.Item(14,4).Value = "'0.6462"
$ws.Cells.Item(14,5).Value = "  +6.96%  "
$ws.Cells.Item(15,4).Value = "'84.18"
$ws.Cells.Item(15,5).Value = "  +9.18%  "
$ws.Cells.Item(16,4).Value = "'4.688"
$ws.Cells.Item(16,5).Value = "  +5.14%  "
$ws.Cells.Item(17,4).Value = "28.898.53"
$ws.Cells.Item(17,5).Value = "  +8.87%  "
$ws.Cells.Item(18,5).Value = "  +0.01%  "
$ws.Cells.Item(19,4).Value = "'0.000007344"
$ws.Cells.Item(19,5).Value = "  +2.69%  "
$ws.Cells.Item(20,4).Value = "'0.9993"
$ws.Cells.Item(20,5).Value = "  +0.23%  "
$ws.Cells.Item(21,5).Value = "  +8.03%  "
$ws.Cells.Item(22,4).Value = "2.044.79"
$ws.Cells.Item(22,5).Value = "  +5.27%  "
$ws.Cells.Item(23,4).Value = "'4.584"
$ws.Cells.Item(23,5).Value = "  +3.87%  "
$ws.Cells.Item(24,4).Value = "'8.848"
$ws.Cells.Item(24,5).Value = "  +3.43%  "
$ws.Cells.Item(25,4).Value = "'5.366"
$ws.Cells.Item(25,5).Value = "  +6.36%  "
$ws.Cells.Item(26,4).Value = "'142.37"
$ws.Cells.Item(26,5).Value = "  +1.90%  "
$ws.Cells.Item(27,4).Value = "'129.17"
$ws.Cells.Item(27,5).Value = "  +20.98%  "
$ws.Cells.Item(28,4).Value = "'16.41"
$ws.Cells.Item(28,5).Value = "  +7.47%  "
$ws.Cells.Item(29,4).Value = "'1.904"
$ws.Cells.Item(29,5).Value = "  +7.42%  "
$ws.Cells.Item(30,4).Value = "'1.415"
$ws.Cells.Item(30,5).Value = "  +2.81%  "
$ws.Cells.Item(31,4).Value = "'4.149"
$ws.Cells.Item(31,5).Value = "  +3.34%  "
$ws.Cells.Item(32,4).Value = "'0.08359"
$ws.Cells.Item(32,5).Value = "  +5.63%  "
$ws.Cells.Item(33,4).Value = "'3.820"
$ws.Cells.Item(33,5).Value = "  +4.13%  "
$ws.Cells.Item(34,4).Value = "'0.04966"
$ws.Cells.Item(34,5).Value = "  +10.38%  "
$ws.Cells.Item(35,4).Value = "'1.099"
$ws.Cells.Item(36,4).Value = "'0.6739"
$ws.Cells.Item(36,5).Value = "  +9.24%  "
$ws.Cells.Item(37,4).Value = "'2.685"
$ws.Cells.Item(37,5).Value = "  +3.51%  "
$ws.Cells.Item(38,4).Value = "'2.317"
$ws.Cells.Item(38,5).Value = "  +15.42%  "
$ws.Cells.Item(39,4).Value = "'2.749"
$ws.Cells.Item(39,5).Value = "  +12.47%  "
$ws.Cells.Item(40,4).Value = "'0.9535"
$ws.Cells.Item(40,5).Value = "  +2.08%  "
$ws.Cells.Item(41,4).Value = "'6.153"
$ws.Cells.Item(41,5).Value = "  +9.62%  "
$ws.Cells.Item(42,4).Value = "'0.01595"
$ws.Cells.Item(42,5).Value = "  +6.86%  "
$ws.Cells.Item(43,4).Value = "'0.9992"
$ws.Cells.Item(43,5).Value = "  +0.12%  "
$ws.Cells.Item(44,4).Value = "'0.4101"
$ws.Cells.Item(44,5).Value = "  +7.16%  "
$ws.Cells.Item(45,4).Value = "'99.88"
$ws.Cells.Item(45,5).Value = "  +0.08%  "
$ws.Cells.Item(46,4).Value = "'7.167"
$ws.Cells.Item(46,5).Value = "  +5.59%  "
$ws.Cells.Item(47,5).Value = "  +6.14%  "
$ws.Cells.Item(48,4).Value = "'0.05516"
$ws.Cells.Item(48,5).Value = "  +2.93%  "
$ws.Cells.Item(49,4).Value = "'8.182"
$ws.Cells.Item(49,5).Value = "  +4.22%  "
$ws.Cells.Item(50,4).Value = "'31.75"
$ws.Cells.Item(50,5).Value = "  +5.76%  "
$ws.Cells.Item(51,5).Value = "  +8.79%  "
